$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B. Excel's ColumnWidth COM property is offset from the raw
# OOXML character-width unit by 5/6, so subtract that to land on width=77.
$ws.Columns.Item(2).ColumnWidth = 76.16666666666667

# --- Apply formatting (wrap text / text number format) BEFORE writing
# --- values so the cells pick up the existing style indices (s="3" for
# --- column B, s="1" for column C) instead of creating new duplicate xfs.
$ws.Range("B2:B6").NumberFormat = "@"
$ws.Range("B2:B6").WrapText = $true
$ws.Range("C3:C6").WrapText = $true

# --- Write new cell contents in the precise order that reproduces the
# --- author's shared-string allocation order (new unique strings must
# --- appear in this sequence: C2, A3, B3, A4, A5, B2, B4, B5, B6, A6).

# 1) C2 -> "benchmark, test_sc_counts" (new shared string)
$ws.Cells.Item(2, 3).Value = "benchmark, test_sc_counts"

# 2) A3 -> "sc.pheno"
$ws.Cells.Item(3, 1).Value = "sc.pheno"

# 3) B3 -> sc.pheno description
$ws.Cells.Item(3, 2).Value = "#' @param sc.pheno data frame with scRNA-Seq profiles as rows, and pheno entries`n#'  in columns. 'nrow(sc.pheno)' must equal 'ncol(sc.counts)'"

# 4) A4 -> "real.counts"
$ws.Cells.Item(4, 1).Value = "real.counts"

# 5) A5 -> "real.props"
$ws.Cells.Item(5, 1).Value = "real.props"

# 6) B2 -> updated sc.counts description
$ws.Cells.Item(2, 2).Value = "#' @param sc.counts non-negative numeric matrix with features as rows, and `n#' scRNA-Seq profiles as columns. 'ncol(sc.counts)' must equal 'nrow(sc.pheno)'"

# 7) B4 -> real.counts description
$ws.Cells.Item(4, 2).Value = "#' @param real.counts non-negative numeric matrix, with features as rows, and `n#' bulk RNA-Seq profiles as columns. 'ncol(sc.counts)' must equal `n#' 'nrow(real.props)'"

# 8) B5 -> real.props description
$ws.Cells.Item(5, 2).Value = "#' @param real.props non-negative numeric matrix, with cell types as rows, `n#' and bulk RNA-Seq profiles."

# 9) B6 -> benchmark.name description
$ws.Cells.Item(6, 2).Value = "#' @param benchmark.name string"

# 10) A6 -> "benchmark.name"
$ws.Cells.Item(6, 1).Value = "benchmark.name"

# --- Fill the remaining (already-existing-string) cells ---
$ws.Cells.Item(3, 3).Value = "benchmark"
$ws.Cells.Item(3, 4).Value = 1

$ws.Cells.Item(4, 3).Value = "benchmark"
$ws.Cells.Item(4, 4).Value = 1

$ws.Cells.Item(5, 3).Value = "benchmark"
$ws.Cells.Item(5, 4).Value = 1

$ws.Cells.Item(6, 3).Value = "benchmark"

$ws.Cells.Item(2, 4).Value = 1

# --- Row heights (auto-fit wrapped content) ---
$ws.Rows.Item(2).RowHeight = 34
$ws.Rows.Item(3).RowHeight = 34
$ws.Rows.Item(4).RowHeight = 51
$ws.Rows.Item(5).RowHeight = 34
$ws.Rows.Item(6).RowHeight = 17

# Update selection to G6
$ws.Range("G6").Select()
